{"js": "const replacements = [\n  [\"2025-09-13 Saturday\", \"2025-09-14 Sunday\"],\n  [\"246\u00f72=123, 0\", \"509\u00f77=72, 5\"],\n  [\"908\u00f73=302, 2\", \"101\u00f78=12, 5\"],\n  [\"881\u00f76=146, 5\", \"756\u00f78=94, 4\"],\n  [\"506\u00f76=84, 2\", \"504\u00f74=126, 0\"],\n  [\"843\u00f77=120, 3\", \"111\u00f79=12, 3\"],\n  [\"909\u00f79=101, 0\", \"601\u00f77=85, 6\"],\n  [\"198\u00f78=24, 6\", \"389\u00f72=194, 1\"],\n  [\"579\u00f79=64, 3\", \"997\u00f73=332, 1\"],\n  [\"782\u00f75=156, 2\", \"335\u00f78=41, 7\"],\n  [\"105\u00f74=26, 1\", \"198\u00f75=39, 3\"],\n  [\"762\u00f73=254, 0\", \"985\u00f78=123, 1\"],\n  [\"448\u00f79=49, 7\", \"185\u00f76=30, 5\"],\n  [\"216\u00f75=43, 1\", \"751\u00f76=125, 1\"],\n  [\"167\u00f77=23, 6\", \"350\u00f76=58, 2\"],\n  [\"592\u00f74=148, 0\", \"455\u00f76=75, 5\"],\n  [\"331\u00f77=47, 2\", \"609\u00f78=76, 1\"],\n  [\"675\u00f72=337, 1\", \"365\u00f73=121, 2\"],\n  [\"436\u00f76=72, 4\", \"427\u00f76=71, 1\"],\n  [\"946\u00f79=105, 1\", \"871\u00f73=290, 1\"],\n  [\"477\u00f76=79, 3\", \"173\u00f77=24, 5\"],\n  [\"119\u00f76=19, 5\", \"292\u00f73=97, 1\"],\n  [\"185\u00f73=61, 2\", \"502\u00f73=167, 1\"],\n  [\"578\u00f79=64, 2\", \"102\u00f74=25, 2\"],\n  [\"433\u00f73=144, 1\", \"732\u00f76=122, 0\"],\n  [\"653\u00f76=108, 5\", \"392\u00f76=65, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2025-09-13 Saturday\", \"2025-09-14 Sunday\")\n    ,@(\"246\u00f72=123, 0\", \"509\u00f77=72, 5\")\n    ,@(\"908\u00f73=302, 2\", \"101\u00f78=12, 5\")\n    ,@(\"881\u00f76=146, 5\", \"756\u00f78=94, 4\")\n    ,@(\"506\u00f76=84, 2\", \"504\u00f74=126, 0\")\n    ,@(\"843\u00f77=120, 3\", \"111\u00f79=12, 3\")\n    ,@(\"909\u00f79=101, 0\", \"601\u00f77=85, 6\")\n    ,@(\"198\u00f78=24, 6\", \"389\u00f72=194, 1\")\n    ,@(\"579\u00f79=64, 3\", \"997\u00f73=332, 1\")\n    ,@(\"782\u00f75=156, 2\", \"335\u00f78=41, 7\")\n    ,@(\"105\u00f74=26, 1\", \"198\u00f75=39, 3\")\n    ,@(\"762\u00f73=254, 0\", \"985\u00f78=123, 1\")\n    ,@(\"448\u00f79=49, 7\", \"185\u00f76=30, 5\")\n    ,@(\"216\u00f75=43, 1\", \"751\u00f76=125, 1\")\n    ,@(\"167\u00f77=23, 6\", \"350\u00f76=58, 2\")\n    ,@(\"592\u00f74=148, 0\", \"455\u00f76=75, 5\")\n    ,@(\"331\u00f77=47, 2\", \"609\u00f78=76, 1\")\n    ,@(\"675\u00f72=337, 1\", \"365\u00f73=121, 2\")\n    ,@(\"436\u00f76=72, 4\", \"427\u00f76=71, 1\")\n    ,@(\"946\u00f79=105, 1\", \"871\u00f73=290, 1\")\n    ,@(\"477\u00f76=79, 3\", \"173\u00f77=24, 5\")\n    ,@(\"119\u00f76=19, 5\", \"292\u00f73=97, 1\")\n    ,@(\"185\u00f73=61, 2\", \"502\u00f73=167, 1\")\n    ,@(\"578\u00f79=64, 2\", \"102\u00f74=25, 2\")\n    ,@(\"433\u00f73=144, 1\", \"732\u00f76=122, 0\")\n    ,@(\"653\u00f76=108, 5\", \"392\u00f76=65, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
